$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = 52.47848103381103

for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 14).Value = $newValue
}
